# Updated cryptos list (Price + Volume(1h) columns) to reflect the latest
# GitHub Actions scrape. Prices/volume % are stored as plain text in the
# sheet (inlineStr), so numeric-looking prices get a leading "'" quote
# prefix to force text entry and avoid Excel coercing them into floats
# (which would mangle trailing zeros / introduce FP noise).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.560.37'
$ws.Range('E2').Value = '  +2.52%  '

$ws.Range('D3').Value = '1.988.24'
$ws.Range('E3').Value = '  +6.03%  '

$ws.Range('D4').Value = "'1.007"
$ws.Range('E4').Value = '  +0.37%  '

$ws.Range('D5').Value = "'329.70"
$ws.Range('E5').Value = '  +1.39%  '

$ws.Range('D6').Value = "'1.006"
$ws.Range('E6').Value = '  +0.31%  '

$ws.Range('D7').Value = "'0.4685"
$ws.Range('E7').Value = '  +1.98%  '

$ws.Range('D8').Value = "'0.3948"
$ws.Range('E8').Value = '  +1.91%  '

$ws.Range('D9').Value = "'46.63"
$ws.Range('E9').Value = '  +0.22%  '

$ws.Range('D10').Value = "'0.07976"
$ws.Range('E10').Value = '  +1.32%  '

$ws.Range('E11').Value = '  +2.03%  '

$ws.Range('D12').Value = "'22.80"
$ws.Range('E12').Value = '  +5.06%  '

$ws.Range('D13').Value = '2.014.42'

$ws.Range('D14').Value = "'7.275"
$ws.Range('E14').Value = '  +4.15%  '

$ws.Range('D15').Value = "'5.890"
$ws.Range('E15').Value = '  +4.32%  '

$ws.Range('D16').Value = "'0.07170"
$ws.Range('E16').Value = '  +3.01%  '

$ws.Range('D17').Value = "'89.03"
$ws.Range('E17').Value = '  +1.07%  '

$ws.Range('E18').Value = '  +0.46%  '

$ws.Range('D19').Value = "'0.000009978"
$ws.Range('E19').Value = '  +0.17%  '

$ws.Range('D20').Value = "'17.38"
$ws.Range('E20').Value = '  +2.65%  '

$ws.Range('D21').Value = "'1.004"
$ws.Range('E21').Value = '  +0.28%  '

$ws.Range('D22').Value = '29.629.13'
$ws.Range('E22').Value = '  +2.71%  '

$ws.Range('D23').Value = "'5.554"
$ws.Range('E23').Value = '  +5.95%  '

$ws.Range('D24').Value = "'11.30"
$ws.Range('E24').Value = '  +3.31%  '

$ws.Range('D25').Value = '2.253.52'
$ws.Range('E25').Value = '  +7.90%  '

$ws.Range('D26').Value = "'2.125"
$ws.Range('E26').Value = '  +1.71%  '

$ws.Range('D27').Value = "'158.27"
$ws.Range('E27').Value = '  +1.65%  '

$ws.Range('D28').Value = "'19.70"
$ws.Range('E28').Value = '  +2.10%  '

$ws.Range('D29').Value = "'5.996"
$ws.Range('E29').Value = '  +0.06%  '

$ws.Range('D30').Value = "'120.61"
$ws.Range('E30').Value = '  +2.91%  '

$ws.Range('D31').Value = "'1.967"
$ws.Range('E31').Value = '  +2.22%  '

$ws.Range('D32').Value = "'0.09453"
$ws.Range('E32').Value = '  +1.32%  '

$ws.Range('D33').Value = "'0.8934"
$ws.Range('E33').Value = '  -0.95%  '

$ws.Range('D34').Value = "'5.294"
$ws.Range('E34').Value = '  +0.75%  '

$ws.Range('D35').Value = "'1.348"
$ws.Range('E35').Value = '  +2.59%  '

$ws.Range('D36').Value = "'3.192"
$ws.Range('E36').Value = '  -2.31%  '

$ws.Range('D37').Value = "'0.05853"
$ws.Range('E37').Value = '  +1.66%  '

$ws.Range('D38').Value = "'1.179"
$ws.Range('E38').Value = '  -0.45%  '

$ws.Range('D39').Value = "'0.02136"
$ws.Range('E39').Value = '  +3.33%  '

$ws.Range('D40').Value = "'7.928"
$ws.Range('E40').Value = '  +3.66%  '

$ws.Range('D41').Value = "'0.5766"
$ws.Range('E41').Value = '  +2.22%  '

$ws.Range('D42').Value = "'0.1826"
$ws.Range('E42').Value = '  +3.64%  '

$ws.Range('D43').Value = "'0.000003104"
$ws.Range('E43').Value = '  +94.47%  '

$ws.Range('D44').Value = "'9.840"
$ws.Range('E44').Value = '  +2.23%  '

$ws.Range('D45').Value = "'12.14"
$ws.Range('E45').Value = '  +2.88%  '

$ws.Range('D46').Value = "'0.5389"
$ws.Range('E46').Value = '  +0.93%  '

$ws.Range('D47').Value = "'2.167"
$ws.Range('E47').Value = '  -3.56%  '

$ws.Range('D48').Value = "'2.648"
$ws.Range('E48').Value = '  +5.61%  '

$ws.Range('D49').Value = "'0.06971"

$ws.Range('D50').Value = "'1.872"
$ws.Range('E50').Value = '  +1.72%  '

$ws.Range('D51').Value = "'114.80"
$ws.Range('E51').Value = '  +1.60%  '
